$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text cell value without letting COM promote
# numeric-looking strings (e.g. "568.83") to real numbers or stamp an
# implicit "Text" number-format style onto the cell. A leading apostrophe
# forces text entry (Excel strips the apostrophe itself); ClearFormats()
# then removes the auto-applied @ style so the cell stays completely
# unstyled, matching the source workbook.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'$value" 
    $range.ClearFormats()
}

# --- Rows 2-51: updated Price / Volume(1h) figures (ranking order unchanged) ---
Set-TextValue $ws.Range('D2') '66.382.62'
Set-TextValue $ws.Range('E2') '  +2.20%  '
Set-TextValue $ws.Range('D3') '3.422.49'
Set-TextValue $ws.Range('E3') '  +1.21%  '
Set-TextValue $ws.Range('E4') '  +0.01%  '
Set-TextValue $ws.Range('D5') '568.83'
Set-TextValue $ws.Range('E5') '  +1.57%  '
Set-TextValue $ws.Range('D6') '181.76'
Set-TextValue $ws.Range('E6') '  +4.74%  '
Set-TextValue $ws.Range('D7') '0.632'
Set-TextValue $ws.Range('E7') '  +1.44%  '
Set-TextValue $ws.Range('D8') '3.417.44'
Set-TextValue $ws.Range('E8') '  +1.40%  '
Set-TextValue $ws.Range('E9') '  -0.09%  '
Set-TextValue $ws.Range('D10') '0.179'
Set-TextValue $ws.Range('E10') '  +6.36%  '
Set-TextValue $ws.Range('D11') '0.643'
Set-TextValue $ws.Range('E11') '  +2.00%  '
Set-TextValue $ws.Range('D12') '55.18'
Set-TextValue $ws.Range('E12') '  +1.95%  '
Set-TextValue $ws.Range('E13') '  +0.73%  '
Set-TextValue $ws.Range('D14') '9.35'
Set-TextValue $ws.Range('E14') '  +2.70%  '
Set-TextValue $ws.Range('D15') '3.973.04'
Set-TextValue $ws.Range('E15') '  +1.38%  '
Set-TextValue $ws.Range('D16') '18.38'
Set-TextValue $ws.Range('E16') '  +0.80%  '
Set-TextValue $ws.Range('D19') '66.391.56'
Set-TextValue $ws.Range('E19') '  +2.25%  '
Set-TextValue $ws.Range('D20') '12.01'
Set-TextValue $ws.Range('E20') '  +1.77%  '
Set-TextValue $ws.Range('E21') '  +1.62%  '
Set-TextValue $ws.Range('D22') '466.46'
Set-TextValue $ws.Range('E22') '  -1.13%  '
Set-TextValue $ws.Range('D23') '5.00'
Set-TextValue $ws.Range('E23') '  +1.02%  '
Set-TextValue $ws.Range('D24') '14.61'
Set-TextValue $ws.Range('E24') '  +8.14%  '
Set-TextValue $ws.Range('D25') '4.16'
Set-TextValue $ws.Range('E25') '  +0.67%  '
Set-TextValue $ws.Range('D26') '89.93'
Set-TextValue $ws.Range('E26') '  +3.57%  '
Set-TextValue $ws.Range('E27') '  +1.81%  '
Set-TextValue $ws.Range('D28') '10.85'
Set-TextValue $ws.Range('E28') '  +0.64%  '
Set-TextValue $ws.Range('E29') '  +1.47%  '
Set-TextValue $ws.Range('D30') '31.43'
Set-TextValue $ws.Range('D31') '6.93'
Set-TextValue $ws.Range('E31') '  +3.56%  '
Set-TextValue $ws.Range('D32') '11.59'
Set-TextValue $ws.Range('E32') '  +0.76%  '
Set-TextValue $ws.Range('D33') '586.40'
Set-TextValue $ws.Range('E33') '  +3.00%  '
Set-TextValue $ws.Range('D34') '62.55'
Set-TextValue $ws.Range('E34') '  +2.06%  '
Set-TextValue $ws.Range('D35') '0.109'
Set-TextValue $ws.Range('E35') '  +1.33%  '
Set-TextValue $ws.Range('E36') '  -0.09%  '
Set-TextValue $ws.Range('E37') '  +4.90%  '
Set-TextValue $ws.Range('D38') '3.61'
Set-TextValue $ws.Range('E38') '  -1.70%  '
Set-TextValue $ws.Range('D39') '36.55'
Set-TextValue $ws.Range('E39') '  +2.70%  '
Set-TextValue $ws.Range('E40') '  +4.00%  '
Set-TextValue $ws.Range('D41') '0.0₃0761'
Set-TextValue $ws.Range('E41') '  +1.71%  '
Set-TextValue $ws.Range('D42') '3.131.66'
Set-TextValue $ws.Range('E42') '  +1.48%  '
Set-TextValue $ws.Range('D43') '2.94'
Set-TextValue $ws.Range('E43') '  +2.92%  '
Set-TextValue $ws.Range('E44') '  +2.36%  '
Set-TextValue $ws.Range('E45') '  +2.89%  '
Set-TextValue $ws.Range('E46') '  +0.26%  '
Set-TextValue $ws.Range('D49') '0.999'
Set-TextValue $ws.Range('E49') '  +0.00%  '
Set-TextValue $ws.Range('D50') '140.94'
Set-TextValue $ws.Range('E50') '  +1.00%  '
Set-TextValue $ws.Range('D51') '8.62'
Set-TextValue $ws.Range('E51') '  +4.07%  '

# --- Rows 17/18 swapped ranking: TRON now ranked above WrappedEther ---
Set-TextValue $ws.Range('B17') 'TRON'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D17') '0.120'
Set-TextValue $ws.Range('E17') '  +0.91%  '
Set-TextValue $ws.Range('B18') 'WrappedEther'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D18') '3.414.47'
Set-TextValue $ws.Range('E18') '  +0.74%  '

# --- Rows 47/48 swapped ranking: dogwifhat now ranked above ApeXProtocol ---
Set-TextValue $ws.Range('B47') 'dogwifhat'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D47') '2.71'
Set-TextValue $ws.Range('E47') '  +17.62%  '
Set-TextValue $ws.Range('B48') 'ApeXProtocol'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range('D48') '3.19'
Set-TextValue $ws.Range('E48') '  +1.68%  '
